$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 2: the StudentID placeholder value had been mistakenly written into
# column C (StudentID) instead of column F (the trailing blank "Rejected"
# reason column that every other row already carries). Clear C2 and give F2
# the same empty-string value used throughout the rest of the column.
$ws.Range("C2").ClearContents()
$ws.Range("F2").Value = ""

# Append the two new rows handed back from the Supervisor / FYP Coordinator
# Class rollover.
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "yes"
$ws.Range("D18").Value = "ASMADHUKUMAR"
$ws.Range("E18").Value = "AVAILABLE"
$ws.Range("F18").Value = ""

$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "hello"
$ws.Range("D19").Value = "ASMADHUKUMAR"
$ws.Range("E19").Value = "AVAILABLE"
$ws.Range("F19").Value = ""

# Reflect the user's final selection in the sheet view.
$ws.Range("B16").Select()
